$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 7.611
$ws.Range("K3").Value = 34.5183
$ws.Range("L3").Value = 10.9108
$ws.Range("M3").Value = 191.321
$ws.Range("N3").Value = 150.448
$ws.Range("O3").Value = 29.7
$ws.Range("P3").Value = 2.54
$ws.Range("R3").Value = 2.71

$ws.Range("Q7").Select()
